# [AFG] added final excel sheets for Afghanistan
#
# 1. Clear the stray empty inline-string cell at "ODI Batting"!B2.
# 2. Add two new worksheets at the end of the workbook:
#      - "ODI Batting Extra" (headers + one data row)
#      - "ODI Bowling Extra" (headers + one data row)
#    The new sheets are produced by copying an existing sheet (so they pick
#    up the same sheetPr/pageMargins as the rest of the workbook) and then
#    clearing + repopulating their contents. Header cells reuse the same
#    bold/centered/bordered header style already used on the other sheets.

$wb = $excel.ActiveWorkbook

# --- 1. Clear ODI Batting!B2 (remove the stray empty inline string cell) ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("B2").ClearContents()

# Template sheet used purely as a source for the header style + sheetPr et al.
$template = $wb.Worksheets.Item("Player Info")
$headerStyleSource = $template.Range("A1")

# --- 2. Add "ODI Batting Extra" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)
$battingExtra = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra.Name = "ODI Batting Extra"
$battingExtra.Cells.Clear()

$battingExtraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $battingExtraHeaders.Length; $i++) {
    $battingExtra.Cells.Item(1, $i + 1).Value = $battingExtraHeaders[$i]
}
$headerStyleSource.Copy()
$battingExtra.Range("A1:F1").PasteSpecial(-4122) # xlPasteFormats

# Row 2 data: MATCH_CODE stays text, BATTING_POSITION is a real number,
# NUM_4/NUM_6/PERCENT_RUNS_OF_TOTAL are left blank, MAN_OF_MATCH is text.
# Force A2/F2 to be stored as text (not auto-coerced to numbers) via a
# temporary "@" text format, then drop back to the "Normal" style so no
# stray number-format style sticks to the cell.
$battingExtra.Range("A2").NumberFormat = "@"
$battingExtra.Range("A2").Value = "4530"
$battingExtra.Range("A2").Style = "Normal"
$battingExtra.Range("B2").Value = 9
$battingExtra.Range("F2").NumberFormat = "@"
$battingExtra.Range("F2").Value = "NO"
$battingExtra.Range("F2").Style = "Normal"

# --- 3. Add "ODI Bowling Extra" sheet after "ODI Batting Extra" ---
$template.Copy($null, $battingExtra)
$bowlingExtra = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra.Name = "ODI Bowling Extra"
$bowlingExtra.Cells.Clear()

$bowlingExtraHeaders = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($i = 0; $i -lt $bowlingExtraHeaders.Length; $i++) {
    $bowlingExtra.Cells.Item(1, $i + 1).Value = $bowlingExtraHeaders[$i]
}
$headerStyleSource.Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122) # xlPasteFormats

# Row 2 data: all three values stored as text.
$bowlingExtra.Range("A2:C2").NumberFormat = "@"
$bowlingExtra.Range("A2").Value = "4530"
$bowlingExtra.Range("B2").Value = "0"
$bowlingExtra.Range("C2").Value = "30.00%"
$bowlingExtra.Range("A2:C2").Style = "Normal"
